$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the HKL-type label list in column B -----------------------
# Two new entries ("Holden" and "Rizzie Spiral") were inserted into the
# master list of HKL types right after "Spiral5". Every label that used
# to sit in rows 4-29 therefore now shows up two rows further down
# (rows 6-31), and "Thomas Hex" was renamed to "Matthies Hex" along the
# way. Column A (the run index) and the C:T result columns for those
# rows are untouched by this relabeling.
$labels = @(
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value2 = $labels[$i]
}

# --- Simulation rerun: two more rows of results ------------------------
# Continuing the existing pattern: column A holds the zero-based run
# index, column B the HKL label for that run, and columns C:T hold a 1
# for every case counted in that run.
$newRows = @(
    @{ Row = 30; Index = 28; Label = "Michael-CCHex" },
    @{ Row = 31; Index = 29; Label = "Michael-SNHex" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Copy the bold/bordered/centered formatting used by every other
    # entry in column A down onto the new row.
    $ws.Cells.Item(29, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value2 = $nr.Index
    $ws.Cells.Item($r, 2).Value2 = $nr.Label
    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value2 = 1
    }
}
